# Updated cryptos list with latest price/volume data, plus a handful of
# coins that re-sorted by rank (row identity swaps) since the previous run.
#
# NOTE: several "Price" values are plain numeric-looking strings (e.g.
# "0.998", "53.43"). Excel.Range.Value would auto-coerce those into
# numbers, losing the fixed string formatting used throughout column D.
# To keep them as text (matching the rest of the sheet, t="inlineStr"/"s"
# cells), we assign them with a leading apostrophe (forces text entry)
# and then reset the cell style back to "Normal" so no stray
# quote-prefix style sticks around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$cellRef,
        [string]$value
    )
    $ws.Range($cellRef).Value = $value
}

function Set-NumericLookingText {
    param(
        [string]$cellRef,
        [string]$value
    )
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "69.922.29"
Set-TextValue "E2" "  -0.02%  "

Set-TextValue "D3" "3.538.61"
Set-TextValue "E3" "  +0.99%  "

Set-NumericLookingText "D4" "0.998"
Set-TextValue "E4" "  +0.01%  "

Set-NumericLookingText "D5" "602.39"
Set-TextValue "E5" "  -0.87%  "

Set-NumericLookingText "D6" "194.94"
Set-TextValue "E6" "  -1.17%  "

Set-TextValue "E7" "  -0.63%  "

Set-TextValue "E9" "  -4.21%  "

Set-NumericLookingText "D10" "0.647"
Set-TextValue "E10" "  -1.72%  "

Set-NumericLookingText "D11" "53.43"
Set-TextValue "E11" "  -1.26%  "

Set-TextValue "E12" "  -1.56%  "

Set-NumericLookingText "D13" "9.49"
Set-TextValue "E13" "  -1.33%  "

Set-TextValue "D14" "4.101.96"
Set-TextValue "E14" "  +1.05%  "

Set-NumericLookingText "D15" "593.47"
Set-TextValue "E15" "  -1.53%  "

# Row 16 now holds Uniswap (was Chainlink)
Set-TextValue "B16" "Uniswap"
Set-TextValue "C16" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-NumericLookingText "D16" "12.79"
Set-TextValue "E16" "  +0.69%  "

# Row 17 now holds WrappedBTC (was Uniswap)
Set-TextValue "B17" "WrappedBTC"
Set-TextValue "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "70.003.74"
Set-TextValue "E17" "  +0.08%  "

# Row 18 now holds Chainlink (was WrappedBTC)
Set-TextValue "B18" "Chainlink"
Set-TextValue "C18" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-NumericLookingText "D18" "19.11"
Set-TextValue "E18" "  +0.52%  "

Set-TextValue "D19" "3.547.66"
Set-TextValue "E19" "  +1.83%  "

Set-TextValue "E20" "  +1.87%  "

Set-NumericLookingText "D21" "0.985"
Set-TextValue "E21" "  -0.86%  "

Set-NumericLookingText "D22" "17.74"
Set-TextValue "E22" "  -0.69%  "

Set-NumericLookingText "D23" "103.20"
Set-TextValue "E23" "  -1.45%  "

Set-NumericLookingText "D24" "5.14"
Set-TextValue "E24" "  +0.75%  "

Set-NumericLookingText "D25" "4.62"
Set-TextValue "E25" "  -0.80%  "

Set-NumericLookingText "D26" "3.05"
Set-TextValue "E26" "  -1.39%  "

Set-NumericLookingText "D27" "10.73"
Set-TextValue "E27" "  -2.34%  "

Set-NumericLookingText "D28" "9.52"
Set-TextValue "E28" "  -3.26%  "

Set-TextValue "E29" "  -2.45%  "

# Row 30 now holds dogwifhat (was NEARProtocol)
Set-TextValue "B30" "dogwifhat"
Set-TextValue "C30" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-NumericLookingText "D30" "4.30"
Set-TextValue "E30" "  -6.83%  "

# Row 31 now holds NEARProtocol (was dogwifhat)
Set-TextValue "B31" "NEARProtocol"
Set-TextValue "C31" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-NumericLookingText "D31" "7.02"
Set-TextValue "E31" "  -2.81%  "

Set-NumericLookingText "D32" "12.32"
Set-TextValue "E32" "  -2.93%  "

Set-TextValue "E33" "  -0.57%  "

Set-TextValue "E34" "  -1.08%  "

# Row 35 now holds Maker (was Fetch.AI)
Set-TextValue "B35" "Maker"
Set-TextValue "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D35" "3.835.40"
Set-TextValue "E35" "  +4.04%  "

# Row 36 now holds Fetch.AI (was Maker)
Set-TextValue "B36" "Fetch.AI"
Set-TextValue "C36" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-NumericLookingText "D36" "3.20"
Set-TextValue "E36" "  +5.54%  "

Set-TextValue "E37" "  +2.72%  "

Set-TextValue "E38" "  +0.18%  "

Set-NumericLookingText "D39" "510.29"
Set-TextValue "E39" "  -2.61%  "

Set-NumericLookingText "D40" "0.390"
Set-TextValue "E40" "  -0.49%  "

Set-TextValue "E41" "  +0.02%  "

Set-NumericLookingText "D42" "36.41"
Set-TextValue "E42" "  -1.52%  "

Set-TextValue "E43" "  -2.69%  "

Set-TextValue "E44" "  -2.59%  "

# Row 45 now holds Stellar (was ThetaToken)
Set-TextValue "B45" "Stellar"
Set-TextValue "C45" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-NumericLookingText "D45" "0.139"
Set-TextValue "E45" "  -0.92%  "

# Row 46 now holds ThetaToken (was Stellar)
Set-TextValue "B46" "ThetaToken"
Set-TextValue "C46" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-NumericLookingText "D46" "2.81"
Set-TextValue "E46" "  -1.73%  "

Set-NumericLookingText "D47" "3.30"
Set-TextValue "E47" "  -0.40%  "

Set-TextValue "E48" "  +0.07%  "

Set-NumericLookingText "D49" "8.50"
Set-TextValue "E49" "  -3.08%  "

Set-TextValue "E50" "  +1.66%  "

Set-TextValue "E51" "  +1.88%  "

$wb.Save()
